$wb = $excel.ActiveWorkbook

# --- Rename sheets (task-order timestamps regenerated) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1650291161330702"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911642753966"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911642773993"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911643564024"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911644304056"

# --- Sheet 1: GNG_TO ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "go_stims-1650291161274727.csv"
$ws.Range("B3").Value = "GNG_stims-16502911612992766.csv"
$ws.Range("B4").Value = "go_stims-16502911613002748.csv"
$ws.Range("B5").Value = "GNG_stims-16502911613295648.csv"

# --- Sheet 2: NB_TO ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = "ZB-match_4-16502911617343986.csv"
$ws.Range("B3").Value = "OB-1650291163303401.csv"
$ws.Range("B4").Value = "OB-16502911626953976.csv"
$ws.Range("B5").Value = "ZB-match_7-16502911618404014.csv"
$ws.Range("B6").Value = "TB-16502911641174035.csv"
$ws.Range("B7").Value = "TB-16502911642543972.csv"
$ws.Range("B8").Value = "OB-16502911622433984.csv"
$ws.Range("B9").Value = "TB-16502911640424023.csv"
$ws.Range("B10").Value = "ZB-match_8-16502911618794007.csv"

# --- Sheet 3: RS_TO (data unchanged, only name changed above) ---

# --- Sheet 4: TOL_TO ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = "MM_stims-16502911643063996.csv"
$ws.Range("B3").Value = "ZM_stims-16502911642803988.csv"
$ws.Range("B4").Value = "MM_stims-16502911643384054.csv"
$ws.Range("B5").Value = "ZM_stims-16502911643074043.csv"
$ws.Range("B6").Value = "MM_stims-16502911643534024.csv"
$ws.Range("B7").Value = "ZM_stims-16502911643394017.csv"

# --- Sheet 5: vSAT_TO ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = "vSAT_stims-16502911644023972.csv"
$ws.Range("B3").Value = "SAT_stims-1650291164365401.csv"
$ws.Range("B4").Value = "vSAT_stims-16502911644173994.csv"
$ws.Range("B5").Value = "SAT_stims-16502911643853993.csv"
